$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column (A) number format down to the new rows before writing values,
# so the new cells pick up style index 1 (numFmtId 14, date m/d/yyyy) like the source rows.
$ws.Range("A674").Copy() | Out-Null
$ws.Range("A675:A689").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$data = @(
    ,@(675, 44907, 3609.10009765625, 3629.239990234375, 3597.340087890625, 3607.06005859375, 25407580, 25407580)
    ,@(676, 44908, 3608.77001953125, 3608.969970703125, 3562.97998046875, 3564.9599609375, 20180190, 20180190)
    ,@(677, 44909, 3578.050048828125, 3584.14990234375, 3549.389892578125, 3555.800048828125, 21631760, 21631760)
    ,@(678, 44910, 3552.389892578125, 3592.31005859375, 3548.2099609375, 3589.489990234375, 20078540, 20078540)
    ,@(679, 44911, 3565.719970703125, 3577.2900390625, 3546.7900390625, 3567.14990234375, 18673392, 18673392)
    ,@(680, 44914, 3565.75, 3586.010009765625, 3513.10009765625, 3525.510009765625, 21680292, 21680292)
    ,@(681, 44915, 3514.239990234375, 3514.570068359375, 3456.780029296875, 3467.97998046875, 18025520, 18025520)
    ,@(682, 44916, 3471.60009765625, 3481.989990234375, 3440.239990234375, 3453.199951171875, 14782428, 14782428)
    ,@(683, 44917, 3476.35009765625, 3489.139892578125, 3433.7900390625, 3443.830078125, 19048650, 19048650)
    ,@(684, 44918, 3421.8701171875, 3450.989990234375, 3409.989990234375, 3426.989990234375, 14672771, 14672771)
    ,@(685, 44921, 3430.8798828125, 3486.5, 3430.8798828125, 3480.330078125, 16333860, 16333860)
    ,@(686, 44922, 3487.590087890625, 3517.739990234375, 3476.340087890625, 3513.85009765625, 16493398, 16493398)
    ,@(687, 44923, 3492.9599609375, 3495.7099609375, 3455.699951171875, 3472.489990234375, 17957264, 17957264)
    ,@(688, 44924, 3460.449951171875, 3499.260009765625, 3456.27001953125, 3476.469970703125, 16271264, 16271264)
    ,@(689, 44925, 3494.699951171875, 3507.47998046875, 3473.530029296875, 3473.530029296875, 15288984, 15288984)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
    $ws.Cells.Item($r, 5).Value2 = $row[5]
    $ws.Cells.Item($r, 6).Formula = "=E" + $r + "/1000"
    $ws.Cells.Item($r, 7).Value2 = $row[6]
    $ws.Cells.Item($r, 8).Value2 = $row[7]
}

# Update the sheet view to match the scrolled/selected state after the new rows were appended.
$excel.ActiveWindow.ScrollRow = 675
$ws.Range("A690").Select() | Out-Null

Write-Output "done"
